# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price cells that are purely numeric-looking text (e.g. "1.00", "9.66") are
# written with a leading apostrophe so Excel keeps them as literal text
# instead of coercing them into numbers (matches the source data, which
# stores prices as plain strings, incl. thousand-dot formats like "68.101.78").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.101.78"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.792.07"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'600.30"
$ws.Range("D6").Value = "'164.96"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'6.52"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "'35.65"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "4.426.97"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "3.769.69"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "68.061.25"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "'18.27"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "'460.87"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'9.66"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D24").Value = "'82.95"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'11.97"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "'9.97"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "3.940.59"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -5.02%  "
$ws.Range("D32").Value = "'7.29"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "'29.24"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "'0.0997"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "'0.986"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D43").Value = "'47.47"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "'43.37"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'151.91"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "'389.41"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'26.61"
$ws.Range("E51").Value = "  -0.74%  "
